$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Last status check on: 19.01.2022 08:00"

$ws.Range("B5").Value = 34.9
$ws.Range("C5").Value = 34.5

$ws.Range("D5").Value = "'+0.4"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "'2022-01-19 08:00:14"
$ws.Range("E5").Style = "Normal"
